# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers section
$ws.Range("D3").Value = 78.40000000000001
$ws.Range("C4").Value = 5159
$ws.Range("D4").Value = 90.90000000000001
$ws.Range("C5").Value = 6655

# Good Drivers section
$ws.Range("B15").Value = 449371
$ws.Range("B16").Value = 14968
$ws.Range("B21").Value = 77999
